$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel PasteSpecial constants
$xlPasteValues  = -4163
$xlPasteFormats = -4122

function Copy-FormatAndValue {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteValues)
    $excel.CutCopyMode = 0
}

function Copy-FormatOnly {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# New "NEW TRAIN/TEST SPLIT" column (D/E) added to the BINARY REFITTING table
# ---------------------------------------------------------------------------

# Row 22: header label (reuse style of L1 / "NEW TRAIN/TEST SPLIT" header)
Copy-FormatOnly "L1" "D22"
$ws.Range("D22").Value = "NEW TRAIN/TEST SPLIT"

# Row 23: "accuracy" column label
Copy-FormatOnly "B24" "D23"
$ws.Range("D23").Value = "accuracy"

# Row 24: bootstrap = true
Copy-FormatAndValue "C24" "D24"
Copy-FormatAndValue "C24" "E24"

# Row 25: max_depth
Copy-FormatOnly "C25" "D25"
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 6

# Row 26: max_features
Copy-FormatAndValue "B26" "D26"
Copy-FormatAndValue "G6" "E26"

# Row 27: max_leaf_nodes
Copy-FormatAndValue "B27" "D27"
$ws.Range("E27").Value = 1000

# Row 28: max_samples
Copy-FormatAndValue "B28" "D28"
Copy-FormatAndValue "C28" "E28"

# Row 29: min_impurity_decrease
Copy-FormatOnly "C29" "D29"
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

# Row 30: min_samples_leaf
Copy-FormatAndValue "B30" "D30"
Copy-FormatOnly "C11" "E30"
$ws.Range("E30").Value = 5

# Row 31: min_samples_split
Copy-FormatAndValue "B31" "D31"
$ws.Range("E31").Value = 3

# Row 32: n_estimators
Copy-FormatAndValue "B32" "D32"
Copy-FormatAndValue "C32" "E32"

# Rows 35-40: new accuracy metrics for the new split
$ws.Range("E35").Value = 0.85
$ws.Range("E36").Value = 0.79
$ws.Range("E37").Value = 0.84
$ws.Range("E38").Value = 0.73
$ws.Range("E39").Value = 0.75
$ws.Range("E40").Value = 0.85

# ---------------------------------------------------------------------------
# New "model fitting process" notes section (rows 43-49)
# ---------------------------------------------------------------------------

Copy-FormatOnly "L1" "A43"
$ws.Range("A43").Value = "model fitting process"

# Build the new (non-bold, colored Consolas) font by starting from the bold
# colored Consolas style (L1) and then switching off bold.
Copy-FormatOnly "L1" "A45"
$ws.Range("A45").Font.Bold = $false
$ws.Range("A45").Value = "1. Bootstrap = True. Otherwise can not set max_samples"

Copy-FormatOnly "A45" "A46"
$ws.Range("A46").Value = "2. Check max_depth effect on overfitting. Set max_depth"

Copy-FormatOnly "A45" "A47"
$ws.Range("A47").Value = "3. General fitting for other parameters: max_features, max_leaf_nodes, max_samples, min_impurity_decrease, min_samples_leaf, min_samples_split"

Copy-FormatOnly "A45" "A48"
$ws.Range("A48").Value = "4. Min_impurity_decrease clearly affect performance. Set to 0"

Copy-FormatOnly "A45" "A49"
$ws.Range("A49").Value = "5. GridSearchCV with other parameters"

Copy-FormatOnly "A45" "A44"
$ws.Range("A44").Value = "0.Splitting with general good results"

# ---------------------------------------------------------------------------
# View state tweaks
# ---------------------------------------------------------------------------
$ws.Range("B11").Select()
$excel.ActiveWindow.ScrollColumn = 1
